$wb = $excel.ActiveWorkbook

# xlPasteFormats
$xlPasteFormats = -4122

# Sheets "展览" (Exhibitions) and "全部类型" (All types) receive identical
# updates, matching what the source data feed regenerated:
#   - F2 (想去人数 / "want to go" count) bumps from 520 to 521
#   - a new row 4 is appended with the 丽水·CCAC动漫游戏嘉年华 event
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Bump the "want to go" counter on row 2 from 520 to 521
    $ws.Range("F2").Value = 521

    # Pre-format B4 as text so the date-like string "2024-10-01" is stored
    # literally (like the other 开始时间 cells) instead of being
    # auto-converted to an Excel date serial number.
    $ws.Range("B4").NumberFormat = "@"

    $ws.Range("A4").Value = 3
    $ws.Range("B4").Value = "2024-10-01"
    $ws.Range("C4").Value = "丽水·CCAC动漫游戏嘉年华"
    $ws.Range("D4").Value = "南环西路109号 九城宴会中心"
    $ws.Range("E4").Value = "2024.10.01 09:00-10.01 16:00"
    $ws.Range("F4").Value = 5
    $ws.Range("G4").Value = 29.9
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=90985"
    $ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202408/IpkQ6J8O1724125974478.jpeg"

    # Stamp row 4's formatting from row 3 (bold/bordered numbering cell in
    # column A, plain General everywhere else). Applied last so it wins
    # over the temporary text format above without disturbing the values.
    $ws.Range("A3:I3").Copy()
    $ws.Range("A4").PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false
